# Scheduled market-data refresh: update currentAveragePrice / Leve price /
# profit columns (H:N) for affected leve rows on each job-sheet, matching
# the latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2094.6365
$ws.Range("I40").Value = 1905.25
$ws.Range("J40").Value = 2599.6667
$ws.Range("K40").Value = 1905.25
$ws.Range("L40").Value = 2599.6667
$ws.Range("M40").Value = -1730.25
$ws.Range("N40").Value = -2949.6667
$ws.Range("H41").Value = 2285.7144
$ws.Range("I41").Value = 4500
$ws.Range("J41").Value = 1400
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 1400
$ws.Range("M41").Value = -4060
$ws.Range("N41").Value = -2280
$ws.Range("H76").Value = 2520
$ws.Range("I76").Value = 2240
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 2240
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -1925
$ws.Range("N76").Value = -3430
$ws.Range("H79").Value = 2520
$ws.Range("I79").Value = 2240
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 2240
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -1148
$ws.Range("N79").Value = -4984
$ws.Range("H121").Value = 768.2
$ws.Range("I121").Value = 200
$ws.Range("J121").Value = 831.3333
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 2493.9999
$ws.Range("M121").Value = 1147
$ws.Range("N121").Value = -5987.9999
$ws.Range("H132").Value = 1249.2727
$ws.Range("I132").Value = 1238.1333
$ws.Range("J132").Value = 1273.1428
$ws.Range("K132").Value = 3714.3999
$ws.Range("L132").Value = 3819.4284
$ws.Range("M132").Value = -1184.3999
$ws.Range("N132").Value = -8879.428400000001
$ws.Range("H137").Value = 1381.9231
$ws.Range("I137").Value = 951.1111
$ws.Range("K137").Value = 2853.3333
$ws.Range("M137").Value = -303.3332999999998
$ws.Range("H138").Value = 3154.0227
$ws.Range("I138").Value = 3913.8948
$ws.Range("J138").Value = 2576.52
$ws.Range("K138").Value = 11741.6844
$ws.Range("L138").Value = 7729.559999999999
$ws.Range("M138").Value = -6601.6844
$ws.Range("N138").Value = -18009.56

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3545.662
$ws.Range("I32").Value = 3055.6428
$ws.Range("K32").Value = 3055.6428
$ws.Range("M32").Value = -2768.6428
$ws.Range("H61").Value = 4557.3438
$ws.Range("I61").Value = 4684.2915
$ws.Range("K61").Value = 4684.2915
$ws.Range("M61").Value = -4472.2915
$ws.Range("H63").Value = 800
$ws.Range("I63").Value = 800
$ws.Range("K63").Value = 800
$ws.Range("M63").Value = -114
$ws.Range("H66").Value = 800
$ws.Range("I66").Value = 800
$ws.Range("K66").Value = 4000
$ws.Range("M66").Value = -568
$ws.Range("H74").Value = 940.94446
$ws.Range("I74").Value = 440.46155
$ws.Range("J74").Value = 2242.2
$ws.Range("K74").Value = 440.46155
$ws.Range("L74").Value = 2242.2
$ws.Range("M74").Value = 433.53845
$ws.Range("N74").Value = -3990.2
$ws.Range("H77").Value = 940.94446
$ws.Range("I77").Value = 440.46155
$ws.Range("J77").Value = 2242.2
$ws.Range("K77").Value = 2202.30775
$ws.Range("L77").Value = 11211
$ws.Range("M77").Value = 2165.69225
$ws.Range("N77").Value = -19947
$ws.Range("H97").Value = 1375.625
$ws.Range("I97").Value = 1301.1666
$ws.Range("K97").Value = 1301.1666
$ws.Range("M97").Value = -805.1666
$ws.Range("H132").Value = 1243.1395
$ws.Range("I132").Value = 947.4595
$ws.Range("K132").Value = 2842.3785
$ws.Range("M132").Value = -312.3785000000003
$ws.Range("H136").Value = 4557.3438
$ws.Range("I136").Value = 4684.2915
$ws.Range("K136").Value = 14052.8745
$ws.Range("M136").Value = -11502.8745

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3493.535
$ws.Range("I134").Value = 3919.5715
$ws.Range("K134").Value = 11758.7145
$ws.Range("M134").Value = -9223.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2117.9
$ws.Range("I31").Value = 1856.7142
$ws.Range("J31").Value = 2727.3333
$ws.Range("K31").Value = 1856.7142
$ws.Range("L31").Value = 2727.3333
$ws.Range("M31").Value = -1561.7142
$ws.Range("N31").Value = -3317.3333
$ws.Range("H34").Value = 2117.9
$ws.Range("I34").Value = 1856.7142
$ws.Range("J34").Value = 2727.3333
$ws.Range("K34").Value = 1856.7142
$ws.Range("L34").Value = 2727.3333
$ws.Range("M34").Value = -1654.7142
$ws.Range("N34").Value = -3131.3333
$ws.Range("H132").Value = 1513.08
$ws.Range("I132").Value = 944.2381
$ws.Range("K132").Value = 2832.7143
$ws.Range("M132").Value = -302.7143000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2499.5715
$ws.Range("I69").Value = 2374.75
$ws.Range("J69").Value = 2666
$ws.Range("K69").Value = 7124.25
$ws.Range("L69").Value = 7998
$ws.Range("M69").Value = -6313.25
$ws.Range("N69").Value = -9620
$ws.Range("H72").Value = 2499.5715
$ws.Range("I72").Value = 2374.75
$ws.Range("J72").Value = 2666
$ws.Range("K72").Value = 21372.75
$ws.Range("L72").Value = 23994
$ws.Range("M72").Value = -17316.75
$ws.Range("N72").Value = -32106
$ws.Range("H122").Value = 794.8889
$ws.Range("I122").Value = 683.75
$ws.Range("K122").Value = 6153.75
$ws.Range("M122").Value = -3703.75
$ws.Range("H131").Value = 28660.08
$ws.Range("I131").Value = 707.5
$ws.Range("J131").Value = 33984.383
$ws.Range("K131").Value = 2122.5
$ws.Range("L131").Value = 101953.149
$ws.Range("M131").Value = 2917.5
$ws.Range("N131").Value = -112033.149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2116.3928
$ws.Range("I102").Value = 1660.762
$ws.Range("J102").Value = 3483.2856
$ws.Range("K102").Value = 1660.762
$ws.Range("L102").Value = 3483.2856
$ws.Range("M102").Value = -38.76199999999994
$ws.Range("N102").Value = -6727.2856
$ws.Range("H132").Value = 1204073
$ws.Range("I132").Value = 1480716.9
$ws.Range("J132").Value = 5282.6665
$ws.Range("K132").Value = 4442150.699999999
$ws.Range("L132").Value = 15847.9995
$ws.Range("M132").Value = -4439620.699999999
$ws.Range("N132").Value = -20907.9995
$ws.Range("H135").Value = 53000
$ws.Range("J135").Value = 53000
$ws.Range("L135").Value = 53000
$ws.Range("N135").Value = -63140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5764.1665
$ws.Range("I16").Value = 6106.364
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 6106.364
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -5936.364
$ws.Range("N16").Value = -2340
$ws.Range("H136").Value = 2265.3845
$ws.Range("J136").Value = 2976.25
$ws.Range("L136").Value = 8928.75
$ws.Range("N136").Value = -14028.75
